$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("message")

# Copy the formatting of the last existing data row (58) down onto the new
# row 59 (this reuses the same style indices / row height instead of
# minting new ones), then fill in the new row's own content.
$ws.Range("A58:C58").Copy()
$ws.Range("A59:C59").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(59).RowHeight = 20

$ws.Range("A59").Formula = "=ROW()-2"
$ws.Range("B59").Value = "食糧を<val1>つ手に入れた"
$ws.Range("C59").Value = "green"
